$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format while assigning, so numeric-looking
# strings (e.g. "1.002", "240.22") are preserved verbatim as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.825.13"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.898.85"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "0.7612"
$ws.Range("E5").Value = "  +3.67%  "
$ws.Range("D6").Value = "240.22"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.3060"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "25.45"
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("D10").Value = "0.06830"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "0.07973"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.908.04"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "0.7406"
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("D14").Value = "5.152"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "90.88"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "29.872.64"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "13.88"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "5.927"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "242.27"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "0.000007678"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "6.929"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "166.87"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "9.210"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "18.68"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").Value = "0.1296"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "2.026"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "1.406"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("D30").Value = "1.514"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("D31").Value = "4.249"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").Value = "4.076"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "0.05248"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").Value = "1.249"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "0.7253"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "2.716"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "0.01922"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D39").Value = "6.151"
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("D40").Value = "0.4406"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "71.79"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "0.8310"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "1.878"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("D45").Value = "7.609"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "9.844"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "99.86"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").Value = "2.048.96"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "35.98"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").Value = "1.478"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").Value = "0.05940"
$ws.Range("E51").Value = "  -0.05%  "

# Restore original (default/general) formatting on column D so cell
# styles match the original workbook.
$ws.Range("D2:D51").ClearFormats()
